# R128.pptx update:
#   - refresh the spec textbox on the "R128-S2" slide (slide id 258) so the
#     CPU / MCU / DSP clock callouts read "up to <freq>" instead of "@<freq>"
#   - stamp the deck with the internal COMMONDATA tag

$p = $ppt.ActivePresentation

# --- slide 2 (SlideID 258, "R128-S2" spec textbox) ---------------------
$slide = $p.Slides.Item(2)
$tr = $slide.Shapes.Item(1).TextFrame.TextRange

$replacements = @(
    @("CPU@600MHz", "CPU up to 600MHz"),
    @("MCU@2", "MCU up to 2"),
    @(" LX7@", " LX7 up to ")
)

foreach ($pair in $replacements) {
    $find = $pair[0]
    $replace = $pair[1]
    $full = $tr.Text
    $idx = $full.IndexOf($find)
    if ($idx -ge 0) {
        $sub = $tr.Characters($idx + 1, $find.Length)
        $sub.Text = $replace
    }
}

# --- custom data tag -----------------------------------------------------
$p.Tags.Add("COMMONDATA", "eyJoZGlkIjoiODUyN2MyMGMxOGVjN2NmMTEyYTIyZGMzOWIxYzMwMzMifQ==")
